$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C numeric updates
$ws.Range("C2").Value = 33
$ws.Range("C6").Value = 27
$ws.Range("C7").Value = 24
$ws.Range("C8").Value = 26
$ws.Range("C9").Value = 29
$ws.Range("C10").Value = 28
$ws.Range("C11").Value = 24
$ws.Range("C12").Value = 42
$ws.Range("C13").Value = 32
$ws.Range("C14").Value = 25
$ws.Range("C15").Value = 26
$ws.Range("C16").Value = 28
$ws.Range("C17").Value = 29
$ws.Range("C18").Value = 28

# Column B text (inline string) updates
$ws.Range("B3").Value = "<downward>"
$ws.Range("B5").Value = "<foot>"
$ws.Range("B9").Value = "<this>"
$ws.Range("B10").Value = "<eight>"
$ws.Range("B13").Value = "<vict>"
$ws.Range("B14").Value = "<up>"
